# Update simulation results for the 380 kV case (pl_mw.xlsx, Sheet1, rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: row number, then new values for columns B, C, D, F, G, L, N
$newValues = @{
    2 = @(2.709022005468171, 0.2366069544716822, 0.09898049955918964, 2.273988311549061, 0.002510418103127687, 0.3047683061611224, 1.871501367930094)
    3 = @(2.54626631303438, 0.205676889423188, 0.09922807900584019, 2.214184979228932, 0.002516829980572283, 0.2940180113051838, 1.882505589821776)
    4 = @(2.448092062838327, 0.1866842450516515, 0.09942768211517716, 2.178965009205157, 0.00252097042540529, 0.2876154528835855, 1.889908560650603)
    5 = @(2.408523108077134, 0.178943351715418, 0.0995209097088221, 2.164986340179041, 0.002522709060067019, 0.2850558453345968, 1.893086991438352)
    6 = @(2.401979033839723, 0.1776578824178614, 0.09953710578390229, 2.162687670710355, 0.002523000867512767, 0.2846338034355966, 1.893624506133087)
    7 = @(2.44755665487645, 0.1865798548057853, 0.09942889138532962, 2.178774979003833, 0.002520993665149098, 0.2875807331980695, 1.889950772656704)
    8 = @(2.652536210656024, 0.2259420356594717, 0.09905594016865393, 2.253054869660133, 0.002512586792099778, 0.3010202626010283, 1.875160953849857)
    9 = @(3.068673561465744, 0.3031627234321377, 0.09870582328736788, 2.410776180849922, 0.002497707084130765, 0.3289649385731366, 1.851325762473294)
    10 = @(3.383393931757325, 0.3599814416625122, 0.09868619558327651, 2.534244209090218, 0.002487741834616458, 0.3504920499079702, 1.837018938087724)
    11 = @(3.528597221417954, 0.3858627990746868, 0.09873001025508188, 2.592114496015711, 0.002483415700881023, 0.3605077037169764, 1.831218303133895)
    12 = @(3.583880392190963, 0.3956694197868273, 0.09875428119533325, 2.614277641061335, 0.002481807085738197, 0.3643328750100494, 1.829124521118928)
    13 = @(3.571960842633075, 0.3935571096544663, 0.0987487109462748, 2.6094932777427, 0.0024821522164525, 0.3635076055613951, 1.829570866088019)
    14 = @(3.533139402157872, 0.3866694714364485, 0.09873185258648931, 2.593932860026371, 0.002483282766838258, 0.360821749902783, 1.83104398030612)
    15 = @(3.509399090796364, 0.3824513937924507, 0.09872252926743386, 2.584434188411308, 0.002483979110771788, 0.3591808255067548, 1.831959723024099)
    16 = @(3.373945950204643, 0.3582908071585962, 0.09868440159032588, 2.530496870041446, 0.002488028709593543, 0.3498420230534691, 1.837412350617527)
    17 = @(3.291374777988494, 0.3434785367496147, 0.09867458727446632, 2.497847263332972, 0.002490565921539195, 0.3441703417921502, 1.840939331435152)
    18 = @(3.24407329150705, 0.3349620729011917, 0.09867390102373008, 2.479228231250801, 0.002492044763664189, 0.3409290923815718, 1.843034498180899)
    19 = @(3.228090516521377, 0.3320790578819128, 0.09867451746759315, 2.472951537053405, 0.002492548829914687, 0.3398352470522639, 1.843755282645404)
    20 = @(3.300144787842555, 0.3450549938269205, 0.09867511813865804, 2.5013062660515, 0.002490293813602406, 0.344771930566921, 1.840556984128014)
    21 = @(3.544534075465663, 0.3886923693060567, 0.09873659511757182, 2.598496544914468, 0.002482949894931696, 0.3616097665067173, 1.830608492830208)
    22 = @(3.705995095545177, 0.417246851002119, 0.09882157984226581, 2.663468777072438, 0.002478322654397915, 0.3728036691469185, 1.824706225832998)
    23 = @(3.619659506241646, 0.4020032585596596, 0.09877209035019519, 2.628657645231954, 0.002480776582750729, 0.3668118032429106, 1.827801161604043)
    24 = @(3.296179337947308, 0.3443422798131337, 0.09867486270515258, 2.499741978822613, 0.002490416770879537, 0.3444998916505995, 1.84072963330587)
    25 = @(2.954544905182161, 0.2822623994508149, 0.09875919428661462, 2.366792957001877, 0.002501561761054494, 0.3212322073820815, 1.857214997081826)
}

$cols = @("B", "C", "D", "F", "G", "L", "N")

foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
